# Joint Acc Creation FCY CAO - Sarib Shamim
# Adds a new "Sheet2" (after the existing "Sheet1") that is a variant of the
# Sheet1 layout with an extra "sbpSubSegment" column inserted after
# "sbpSegment", and makes Sheet2 the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- Update the selection on Sheet1 (it loses tabSelected to Sheet2 below) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("F2").Select() | Out-Null

# --- Create the new worksheet after Sheet1 ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Sheet2"

# --- Populate Sheet2, column by column (header row then data row) ---

# A: Cid
$ws2.Range("A1").Value = "Cid"
$ws2.Range("A2").Value = 11241908

# B: CategoryProduct
$ws2.Range("B1").Value = "CategoryProduct"
$ws2.Range("B2").Value = "6-012"

# C: Acc name
$ws2.Range("C1").Value = "Acc name"
$ws2.Range("C2").Value = "abc"

# D: Acc name2
$ws2.Range("D1").Value = "Acc name2"
$ws2.Range("D2").Value = "abc"

# E: SignOffData (numeric value kept, but stored with a Text number format)
$ws2.Range("E1").Value = "SignOffData"
$ws2.Range("E2").Value = 20230106
$ws2.Range("E2").NumberFormat = "@"

# F: sbpCompany
$ws2.Range("F1").Value = "sbpCompany"
$ws2.Range("F2").Value = 4

# G: sbpSector
$ws2.Range("G1").Value = "sbpSector"
$ws2.Range("G2").Value = 3

# H: sbpSubSector (text value that looks like a date -> force text format first)
$ws2.Range("H1").Value = "sbpSubSector"
$ws2.Range("H2").NumberFormat = "@"
$ws2.Range("H2").Value = "3-01"

# I: sbpSegment (text value that looks like a date -> force text format first)
$ws2.Range("I1").Value = "sbpSegment"
$ws2.Range("I2").NumberFormat = "@"
$ws2.Range("I2").Value = "3-01-01"

# J: sbpSubSegment (new column; text value -> force text format first)
$ws2.Range("J1").Value = "sbpSubSegment"
$ws2.Range("J2").NumberFormat = "@"
$ws2.Range("J2").Value = "3-01-01-0100"

# K: Jholder
$ws2.Range("K1").Value = "Jholder"
$ws2.Range("K2").Value = 16206304

# L: Purpose
$ws2.Range("L1").Value = "Purpose"
$ws2.Range("L2").Value = "testing"

# M: ExpectedNumOfTxn
$ws2.Range("M1").Value = "ExpectedNumOfTxn"
$ws2.Range("M2").Value = 20

# N: AC Screen list
$ws2.Range("N1").Value = "AC Screen list"
$ws2.Range("N2").Value = 2

# O: RelationCode
$ws2.Range("O1").Value = "RelationCode"
$ws2.Range("O2").Value = 6

# P: TurnoverM
$ws2.Range("P1").Value = "TurnoverM"
$ws2.Range("P2").Value = "Below 1M"

# Q: TurnoverA
$ws2.Range("Q1").Value = "TurnoverA"
$ws2.Range("Q2").Value = "1M to 5M"

# R: debitTxnNum
$ws2.Range("R1").Value = "debitTxnNum"
$ws2.Range("R2").Value = 20

# S: TurnoverDebitMonth
$ws2.Range("S1").Value = "TurnoverDebitMonth"
$ws2.Range("S2").Value = "Below 1M"

# --- Make Sheet2 the active tab with the appropriate selection ---
$ws2.Activate() | Out-Null
$ws2.Range("J2").Select() | Out-Null
